$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.238101243972778
$ws.Range("B1").Value = 2.299833297729492
$ws.Range("C1").Value = 3.426994323730469
$ws.Range("D1").Value = 3.851266384124756
$ws.Range("E1").Value = 1.038953900337219
